# chore: update Sheets via scheduled runner
# Refresh market-board price/profit figures on each Leve Profits worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 7315.75
$ws.Range("I64").Value = 4733.5713
$ws.Range("J64").Value = 8706.154
$ws.Range("K64").Value = 4733.5713
$ws.Range("L64").Value = 8706.154
$ws.Range("M64").Value = -4485.5713
$ws.Range("N64").Value = -9202.154
# Row 67
$ws.Range("H67").Value = 7315.75
$ws.Range("I67").Value = 4733.5713
$ws.Range("J67").Value = 8706.154
$ws.Range("K67").Value = 4733.5713
$ws.Range("L67").Value = 8706.154
$ws.Range("M67").Value = -3875.5713
$ws.Range("N67").Value = -10422.154
# Row 74
$ws.Range("H74").Value = 7786.864
$ws.Range("I74").Value = 3500.6
$ws.Range("K74").Value = 3500.6
$ws.Range("M74").Value = -2564.6
# Row 77
$ws.Range("H77").Value = 7786.864
$ws.Range("I77").Value = 3500.6
$ws.Range("K77").Value = 17503
$ws.Range("M77").Value = -12823
# Row 98
$ws.Range("H98").Value = 1521.7097
$ws.Range("I98").Value = 1537.5714
$ws.Range("K98").Value = 1537.5714
$ws.Range("M98").Value = -39.57140000000004
# Row 103
$ws.Range("H103").Value = 1219.75
$ws.Range("I103").Value = 959.6667
$ws.Range("K103").Value = 2879.0001
$ws.Range("M103").Value = -2293.0001
# Row 112
$ws.Range("H112").Value = 5374.591
$ws.Range("J112").Value = 5429.5845
$ws.Range("L112").Value = 16288.7535
$ws.Range("N112").Value = -18504.7535
# Row 122
$ws.Range("H122").Value = 1521.7097
$ws.Range("I122").Value = 1537.5714
$ws.Range("K122").Value = 4612.7142
$ws.Range("M122").Value = -2162.7142
# Row 132
$ws.Range("H132").Value = 63657.76
$ws.Range("I132").Value = 66499.72
$ws.Range("K132").Value = 199499.16
$ws.Range("M132").Value = -196969.16
# Row 138
$ws.Range("H138").Value = 2614.7334
$ws.Range("J138").Value = 2803.682
$ws.Range("L138").Value = 8411.045999999998
$ws.Range("N138").Value = -18691.046

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1711.6
$ws.Range("I45").Value = 1701.7778
$ws.Range("K45").Value = 1701.7778
$ws.Range("M45").Value = -1324.7778
# Row 74
$ws.Range("H74").Value = 3265.1428
$ws.Range("I74").Value = 1024.8334
$ws.Range("J74").Value = 7297.7
$ws.Range("K74").Value = 1024.8334
$ws.Range("L74").Value = 7297.7
$ws.Range("M74").Value = -150.8334
$ws.Range("N74").Value = -9045.700000000001
# Row 77
$ws.Range("H77").Value = 3265.1428
$ws.Range("I77").Value = 1024.8334
$ws.Range("J77").Value = 7297.7
$ws.Range("K77").Value = 5124.166999999999
$ws.Range("L77").Value = 36488.5
$ws.Range("M77").Value = -756.1669999999995
$ws.Range("N77").Value = -45224.5
# Row 110
$ws.Range("H110").Value = 2163
$ws.Range("I110").Value = 1914.5
$ws.Range("J110").Value = 3778.25
$ws.Range("K110").Value = 1914.5
$ws.Range("L110").Value = 3778.25
$ws.Range("M110").Value = 130.5
$ws.Range("N110").Value = -7868.25
# Row 132
$ws.Range("H132").Value = 365550.03
$ws.Range("I132").Value = 401605.53
$ws.Range("K132").Value = 1204816.59
$ws.Range("M132").Value = -1202286.59

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4219.696
$ws.Range("I105").Value = 4121.095
$ws.Range("K105").Value = 4121.095
$ws.Range("M105").Value = -2374.095

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11277.4
$ws.Range("I31").Value = 4034.2258
$ws.Range("J31").Value = 27315.857
$ws.Range("K31").Value = 4034.2258
$ws.Range("L31").Value = 27315.857
$ws.Range("M31").Value = -3739.2258
$ws.Range("N31").Value = -27905.857
# Row 34
$ws.Range("H34").Value = 11277.4
$ws.Range("I34").Value = 4034.2258
$ws.Range("J34").Value = 27315.857
$ws.Range("K34").Value = 4034.2258
$ws.Range("L34").Value = 27315.857
$ws.Range("M34").Value = -3832.2258
$ws.Range("N34").Value = -27719.857
# Row 122
$ws.Range("H122").Value = 2733.8235
$ws.Range("I122").Value = 2540.65
$ws.Range("J122").Value = 3009.7856
$ws.Range("K122").Value = 7621.950000000001
$ws.Range("L122").Value = 9029.356800000001
$ws.Range("M122").Value = -5171.950000000001
$ws.Range("N122").Value = -13929.3568
# Row 132
$ws.Range("H132").Value = 35900504
$ws.Range("I132").Value = 47621400
$ws.Range("J132").Value = 737805.3
$ws.Range("K132").Value = 142864200
$ws.Range("L132").Value = 2213415.9
$ws.Range("M132").Value = -142861670
$ws.Range("N132").Value = -2218475.9

$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 4990
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 4990
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 14970
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -16592
# Row 72
$ws.Range("H72").Value = 4990
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 4990
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 44910
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -53022
# Row 137
$ws.Range("H137").Value = 1983.8462
$ws.Range("J137").Value = 3746.25
$ws.Range("L137").Value = 11238.75
$ws.Range("N137").Value = -21438.75

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 834.8929000000001
$ws.Range("I97").Value = 846.2917
$ws.Range("J97").Value = 766.5
$ws.Range("K97").Value = 846.2917
$ws.Range("L97").Value = 766.5
$ws.Range("M97").Value = -350.2917
$ws.Range("N97").Value = -1758.5
# Row 122
$ws.Range("H122").Value = 30453.18
$ws.Range("I122").Value = 43095.8
$ws.Range("J122").Value = 7877.0713
$ws.Range("K122").Value = 129287.4
$ws.Range("L122").Value = 23631.2139
$ws.Range("M122").Value = -126837.4
$ws.Range("N122").Value = -28531.2139
# Row 126
$ws.Range("H126").Value = 1391743
$ws.Range("I126").Value = 2383369.8
$ws.Range("J126").Value = 3465.6
$ws.Range("K126").Value = 7150109.399999999
$ws.Range("L126").Value = 10396.8
$ws.Range("M126").Value = -7147639.399999999
$ws.Range("N126").Value = -15336.8
# Row 132
$ws.Range("H132").Value = 28921122
$ws.Range("I132").Value = 34902500
$ws.Range("J132").Value = 11132.667
$ws.Range("K132").Value = 104707500
$ws.Range("L132").Value = 33398.001
$ws.Range("M132").Value = -104704970
$ws.Range("N132").Value = -38458.001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3695.3462
$ws.Range("I7").Value = 3403.8
$ws.Range("K7").Value = 3403.8
$ws.Range("M7").Value = -3291.8
# Row 40
$ws.Range("H40").Value = 2856.6667
$ws.Range("I40").Value = 2837
$ws.Range("K40").Value = 2837
$ws.Range("M40").Value = -2701
# Row 74
$ws.Range("H74").Value = 76173.60000000001
$ws.Range("J74").Value = 76173.60000000001
$ws.Range("L74").Value = 76173.60000000001
$ws.Range("N74").Value = -78169.60000000001
# Row 77
$ws.Range("H77").Value = 76173.60000000001
$ws.Range("J77").Value = 76173.60000000001
$ws.Range("L77").Value = 228520.8
$ws.Range("N77").Value = -238504.8
# Row 122
$ws.Range("H122").Value = 5257.6313
$ws.Range("I122").Value = 4793.1333
$ws.Range("J122").Value = 6999.5
$ws.Range("K122").Value = 14379.3999
$ws.Range("L122").Value = 20998.5
$ws.Range("M122").Value = -11929.3999
$ws.Range("N122").Value = -25898.5
# Row 126
$ws.Range("H126").Value = 3695.3462
$ws.Range("I126").Value = 3403.8
$ws.Range("K126").Value = 10211.4
$ws.Range("M126").Value = -7741.400000000001
# Row 132
$ws.Range("H132").Value = 1745454.8
$ws.Range("I132").Value = 2179742.5
$ws.Range("K132").Value = 6539227.5
$ws.Range("M132").Value = -6536697.5
# Row 136
$ws.Range("H136").Value = 2822.111
$ws.Range("I136").Value = 1535.7059
$ws.Range("K136").Value = 4607.1177
$ws.Range("M136").Value = -2057.1177

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 1015001
$ws.Range("J5").Value = 30002
$ws.Range("L5").Value = 30002
$ws.Range("N5").Value = -30226
# Row 122
$ws.Range("H122").Value = 2542.0527
$ws.Range("I122").Value = 2169
$ws.Range("K122").Value = 6507
$ws.Range("M122").Value = -4057
# Row 126
$ws.Range("H126").Value = 5590.231
$ws.Range("I126").Value = 5267.3
$ws.Range("K126").Value = 15801.9
$ws.Range("M126").Value = -13331.9
# Row 132
$ws.Range("H132").Value = 6291696
$ws.Range("I132").Value = 6941861
$ws.Range("K132").Value = 20825583
$ws.Range("M132").Value = -20823053
